$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("multiple")

# The "multiple" sheet originally held two stacked tables starting at row 6
# (Entrate block) and row 31 (Uscite block), leaving rows 1-5 blank above
# the first table. Remove those 5 leading blank rows so the data starts at
# row 1 (shifting everything - including the merged cell ranges - up by 5).
$ws.Rows("1:5").Delete()
